$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert new row 4: dataset.preview.table ---
$ws.Range("A4:B4").EntireRow.Insert() | Out-Null
$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"
$ws.Rows.Item(4).RowHeight = 120
$ws.Range("A4:B4").WrapText = $true
$ws.Range("A4:B4").VerticalAlignment = -4108

# --- Insert new row 5: dataset.preview.line ---
$ws.Range("A5:B5").EntireRow.Insert() | Out-Null
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"
$ws.Rows.Item(5).RowHeight = 120
$ws.Range("A5:B5").WrapText = $true
$ws.Range("A5:B5").VerticalAlignment = -4108

# --- Update the view: scroll/selection on B10, top-left cell reset ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("B10").Select() | Out-Null
